# Update Work Breakdown Agreement workbook for Assignment 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Column width tweaks (minor re-flow from the resave)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5
$ws.Columns.Item(2).ColumnWidth = 31.83
$ws.Columns.Item(3).ColumnWidth = 21.83
$ws.Columns.Item(4).ColumnWidth = 21.5
$ws.Columns.Item(7).ColumnWidth = 23

# ---------------------------------------------------------------
# Row height tweaks on existing rows
# ---------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28
$ws.Rows.Item(2).RowHeight = 13
$ws.Rows.Item(3).RowHeight = 13
$ws.Rows.Item(4).RowHeight = 13
$ws.Rows.Item(5).RowHeight = 28
$ws.Rows.Item(6).RowHeight = 13
$ws.Rows.Item(7).RowHeight = 13
$ws.Rows.Item(8).RowHeight = 13
$ws.Rows.Item(9).RowHeight = 13
$ws.Rows.Item(12).RowHeight = 13
$ws.Rows.Item(13).RowHeight = 26
$ws.Rows.Item(20).RowHeight = 28

# ---------------------------------------------------------------
# New "ASSIGNMENT 3" section
# ---------------------------------------------------------------

# Row 24 - section header + first task
$ws.Range("A24").Value = "ASSIGNMENT 3"
$ws.Range("A24").Font.Name = "Arial"
$ws.Range("A24").Font.Size = 10
$ws.Range("A24").Font.Bold = $true
$ws.Range("A24").Font.Color = 0

$ws.Range("B24").Value = "Add new map"
$ws.Range("C24").Value = "Harun"
$ws.Range("D24").Value = "Sara"

# Row 25
$ws.Range("B25").Value = "Add Water (Ground)"
$ws.Range("C25").Value = "Harun"
$ws.Range("D25").Value = "Sara"

# Row 26
$ws.Range("B26").Value = "Add Reeds growing and Fish spawning"
$ws.Range("C26").Value = "Sara"
$ws.Range("D26").Value = "Harun"

# Row 29 values are entered here (matches original authoring order so that
# new shared-string entries come out in the same order as the source file)
$ws.Range("B29").Value = "T-Rex"
$ws.Range("C29").Value = "Harun"
$ws.Range("C29").WrapText = $true
$ws.Range("D29").Value = "Sara"

# Row 30
$ws.Range("B30").Value = "Quit and Gameover"
$ws.Range("C30").Value = "Sara"
$ws.Range("D30").Value = "Harun"

# Row 28 - wrapped multi-line description (Pteanodons)
$ws.Range("B28").Value = "Pteanodons:`n- Flying over Land and Water`n- Move up to 2 squares per turn --> need to modify WanderBehaviour"
$ws.Range("B28").WrapText = $true
$ws.Range("C28").Value = "Harun"
$ws.Range("D28").Value = "Sara"
$ws.Rows.Item(28).RowHeight = 56

# Row 27 - wrapped multi-line description (Plesiosaurs)
$ws.Range("B27").Value = "Plesiosaurs:`n- marine egg hatching`n- food preference --> need to modify SeekFoodBehaviour"
$ws.Range("B27").WrapText = $true
$ws.Range("C27").Value = "Sara"
$ws.Range("D27").Value = "Harun"
$ws.Rows.Item(27).RowHeight = 56

# Row 31 - sign-off prompt (second copy, like row 20/21 & row 5/6)
$ws.Range("A31").Value = "Do you accept this WBA? [Sign Below]"
$ws.Range("A31").Font.Name = "Arial"
$ws.Range("A31").Font.Size = 10
$ws.Range("A31").Font.Bold = $true
$ws.Range("A31").WrapText = $true
$ws.Range("B31").Font.Name = "Arial"
$ws.Range("B31").Font.Size = 10
$ws.Rows.Item(31).RowHeight = 28

# Row 32 - signature
$ws.Range("A32").Value = "Yes - Sara"
$ws.Rows.Item(32).RowHeight = 15.75

# ---------------------------------------------------------------
# Restore the active selection similar to the authored file
# ---------------------------------------------------------------
$ws.Range("B36").Select()
